$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style (s="2") from the last existing data row (A1084)
# down through the new rows so the new timestamp cells pick up the same
# numFmt/font/border/alignment (YYYY-MM-DD HH:MM:SS) without creating a new style.
$ws.Range("A1084").Copy() | Out-Null
$ws.Range("A1085:A1116").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$data = @(
    @(1085, 45534.5, 1.435, 1.455, 1.325, 1.335, 1312671.85),
    @(1086, 45534.66666666666, 1.335, 1.408, 1.324, 1.366, 679585.28),
    @(1087, 45534.83333333334, 1.366, 1.427, 1.364, 1.426, 490037.79),
    @(1088, 45535, 1.426, 1.509, 1.414, 1.467, 1176305.85),
    @(1089, 45535.16666666666, 1.467, 1.49, 1.457, 1.475, 451056.44),
    @(1090, 45535.33333333334, 1.476, 1.487, 1.45, 1.465, 470529.31),
    @(1091, 45535.5, 1.465, 1.465, 1.366, 1.371, 1283913.99),
    @(1092, 45535.66666666666, 1.372, 1.372, 1.337, 1.352, 444145.47),
    @(1093, 45535.83333333334, 1.352, 1.368, 1.343, 1.361, 140031.81),
    @(1094, 45536, 1.362, 1.363, 1.336, 1.338, 367752.88),
    @(1095, 45536.16666666666, 1.338, 1.349, 1.324, 1.335, 437562.48),
    @(1096, 45536.33333333334, 1.334, 1.335, 1.304, 1.318, 534543.2),
    @(1097, 45536.5, 1.317, 1.324, 1.276, 1.319, 618235.34),
    @(1098, 45536.66666666666, 1.319, 1.371, 1.309, 1.334, 1114321.96),
    @(1099, 45536.83333333334, 1.334, 1.346, 1.257, 1.287, 1106161.51),
    @(1100, 45537, 1.286, 1.297, 1.261, 1.281, 438244.37),
    @(1101, 45537.16666666666, 1.281, 1.282, 1.243, 1.246, 928405.01),
    @(1102, 45537.33333333334, 1.245, 1.293, 1.238, 1.27, 1005552.52),
    @(1103, 45537.5, 1.269, 1.278, 1.236, 1.249, 451272.29),
    @(1104, 45537.66666666666, 1.249, 1.266, 1.236, 1.248, 479026.58),
    @(1105, 45537.83333333334, 1.248, 1.257, 1.246, 1.255, 18111.08),
    @(1106, 45538, 1.268, 1.287, 1.254, 1.257, 498418.6),
    @(1107, 45538.16666666666, 1.258, 1.264, 1.244, 1.253, 219458.47),
    @(1108, 45538.33333333334, 1.253, 1.254, 1.222, 1.225, 510235.11),
    @(1109, 45538.5, 1.225, 1.236, 1.171, 1.174, 1065157.81),
    @(1110, 45538.66666666666, 1.174, 1.2, 1.17, 1.188, 644460.8),
    @(1111, 45538.83333333334, 1.187, 1.193, 1.169, 1.172, 364023.89),
    @(1112, 45539, 1.172, 1.187, 1.121, 1.183, 1095443.53),
    @(1113, 45539.16666666666, 1.183, 1.209, 1.177, 1.207, 785335.5699999999),
    @(1114, 45539.33333333334, 1.208, 1.208, 1.2, 1.208, 35575.53),
    @(1115, 45539.5, 1.199, 1.24, 1.196, 1.239, 639428.25),
    @(1116, 45539.66666666666, 1.239, 1.311, 1.233, 1.238, 1757333.16)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
